$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'21.691.68"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'1.533.64"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "'288.63"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'0.3927"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("D8").Value = "'0.3156"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "'42.03"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").Value = "'0.07167"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").Value = "'1.044"
$ws.Range("E11").Value = "  -7.09%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'5.623"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "'18.46"
$ws.Range("D15").Value = "'6.600"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "'1.537.28"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "'0.00001096"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'0.06593"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'82.95"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'6.107"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").Value = "'10.83"
$ws.Range("E23").Value = "  -5.67%  "
$ws.Range("D24").Value = "'2.378"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").Value = "'21.699.34"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value = "'2.345"
$ws.Range("D27").Value = "'146.57"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").Value = "'18.31"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").Value = "'4.840"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "'1.709.88"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "'117.19"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.877"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9579"
$ws.Range("E33").Value = "  -13.64%  "
$ws.Range("D34").Value = "'0.08164"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'8.623"
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("E36").Value = "  -2.23%  "
$ws.Range("D37").Value = "'5.096"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").Value = "'0.02193"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").Value = "'1.438"
$ws.Range("E39").Value = "  -12.55%  "
$ws.Range("D40").Value = "'0.2021"
$ws.Range("E40").Value = "  -4.18%  "
$ws.Range("D41").Value = "'1.175"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'10.73"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "'0.5710"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.99"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.728"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "'0.5474"
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("D48").Value = "'1.155"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "'115.81"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").Value = "'1.859"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").Value = "'0.06695"
$ws.Range("E51").Value = "  -2.68%  "
